$d = $word.ActiveDocument

# All replacements listed in document order. Using wdReplaceOne (1) and
# always searching from the very start of the document (Find.Execute here
# is effectively document-global) guarantees that repeated/duplicate source
# strings (e.g. "29÷3=") get matched and replaced one at a time, in the
# correct left-to-right order, without clobbering later occurrences that
# must become a different value.
$replacements = @(
    @{ Old = "2026-02-14 Saturday"; New = "2026-02-15 Sunday" },

    @{ Old = "72÷4="; New = "91÷6=" },
    @{ Old = "75÷5="; New = "51÷2=" },
    @{ Old = "84÷7="; New = "42÷7=" },
    @{ Old = "54÷9="; New = "34÷8=" },
    @{ Old = "42÷2="; New = "32÷8=" },

    @{ Old = "48÷5="; New = "88÷9=" },
    @{ Old = "94÷4="; New = "61÷8=" },
    @{ Old = "62÷6="; New = "44÷7=" },
    @{ Old = "62÷8="; New = "49÷5=" },
    @{ Old = "48÷2="; New = "76÷7=" },

    @{ Old = "80÷8="; New = "53÷4=" },
    @{ Old = "29÷3="; New = "84÷5=" },
    @{ Old = "32÷5="; New = "56÷2=" },
    @{ Old = "93÷4="; New = "15÷7=" },
    @{ Old = "88÷3="; New = "19÷7=" },

    @{ Old = "63÷4="; New = "21÷9=" },
    @{ Old = "66÷2="; New = "24÷2=" },
    @{ Old = "20÷5="; New = "75÷7=" },
    @{ Old = "29÷3="; New = "66÷5=" },
    @{ Old = "32÷4="; New = "97÷8=" },

    @{ Old = "33÷9="; New = "72÷5=" },
    @{ Old = "41÷4="; New = "87÷3=" },
    @{ Old = "17÷5="; New = "49÷8=" },
    @{ Old = "22÷6="; New = "10÷5=" },
    @{ Old = "66÷3="; New = "25÷5=" }
)

foreach ($rep in $replacements) {
    $d.Content.Find.Execute($rep.Old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $rep.New, 1)
}
